$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 4 days,
# keeping the same time-of-day fraction (re-training window moved from
# 2025-09-19 to 2025-09-23).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    $cell.Value = $current + 4
}

# Updated solar production readings (column B) for the 26-42 window.
$newB = @{
    26 = 0
    27 = 13
    28 = 58
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
}

foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}
